# Applies the row-content permutation described by the commit diff.
# For each affected group of adjacent rows, the values in columns B..AD
# (everything except the running index in column A) are permuted among
# the rows of the group, while column A (the row index) stays put.
#
# Groups (destination row -> source row whose B:AD content it receives):
#   20 <- 21, 21 <- 20                     (simple swap)
#   34 <- 35, 35 <- 34                     (simple swap)
#   40 <- 42, 42 <- 40                     (simple swap)
#   99 <- 100, 100 <- 99                   (simple swap)
#   116 <- 118, 117 <- 116, 118 <- 117     (3-way rotation)
#   128 <- 129, 129 <- 128                 (simple swap)
#   166 <- 167, 167 <- 166                 (simple swap)
#   176 <- 177, 177 <- 176                 (simple swap)
#   188 <- 189, 189 <- 190, 190 <- 188     (3-way rotation)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Denmark Division 1")

# Columns B (2) through AD (30)
$firstCol = 2
$lastCol = 30

# Capture "before" snapshots (cell by cell, so each value keeps its exact
# original representation) for every row that participates in a permutation,
# before any writes happen.
$rowsToCapture = @(20,21,34,35,40,42,99,100,116,117,118,128,129,166,167,176,177,188,189,190)
$snapshot = @{}
foreach ($r in $rowsToCapture) {
    $rowVals = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# destination -> source mapping
$mapping = @{
    20  = 21
    21  = 20
    34  = 35
    35  = 34
    40  = 42
    42  = 40
    99  = 100
    100 = 99
    116 = 118
    117 = 116
    118 = 117
    128 = 129
    129 = 128
    166 = 167
    167 = 166
    176 = 177
    177 = 176
    188 = 189
    189 = 190
    190 = 188
}

foreach ($dest in $mapping.Keys) {
    $src = $mapping[$dest]
    $srcVals = $snapshot[$src]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $newVal = $srcVals[$c]
        $cell = $ws.Cells.Item($dest, $c)
        # Only touch cells whose value actually changes, to avoid needlessly
        # rewriting (and thus reformatting) cells that already hold the
        # target value.
        if ($cell.Value2 -ne $newVal) {
            $cell.Value2 = $newVal
        }
    }
}

$wb.Save()
